$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
ECs	Efemp1	Egfr	ECs	3	1	2.325466	6.976398000000001	0.03804623361484849	0.0380462336148485	3	1	1.701929666666667	5.105789	0.02105622887134972	0.02105622887134972	3.957779574224666	35.620016168022	0.000801110202687089	0.0008011102026870891
ECs	Efemp1	Egfr	FAPs	3	1	2.325466	6.976398000000001	0.03804623361484849	0.0380462336148485	3	1	62.503947	187.511841	0.7732971809418951	0.7732971809418953	145.350803614302	1308.157232528718	0.02942104519981911	0.02942104519981912
ECs	Efemp1	Egfr	Inflammatory-Mac	3	1	2.325466	6.976398000000001	0.03804623361484849	0.0380462336148485	3	1	0.4338690000000001	1.301607	0.005367815805265532	0.005367815805265533	1.008947607954	9.080528471586002	0.0002042251741286085	0.0002042251741286086
ECs	Efemp1	Egfr	MuSCs	3	1	2.325466	6.976398000000001	0.03804623361484849	0.0380462336148485	3	1	15.972384	47.917152	0.1976099128607259	0.1976099128607259	37.143235930944	334.289123378496	0.007518312909309031	0.007518312909309034
ECs	Efemp1	Egfr	Resolving-Mac	3	1	2.325466	6.976398000000001	0.03804623361484849	0.0380462336148485	3	1	0.2157183333333333	0.647155	0.002668861520763652	0.002668861520763652	0.5016456497433334	4.514810847690001	0.0001015401289046537	0.0001015401289046537
FAPs	Efemp1	Egfr	ECs	3	1	58.41791133333334	175.253734	0.9557574704938732	0.9557574704938733	3	1	1.701929666666667	5.105789	0.02105622887134972	0.02105622887134972	99.42317636290289	894.808587266126	0.02012464804422127	0.02012464804422127
FAPs	Efemp1	Egfr	FAPs	3	1	58.41791133333334	175.253734	0.9557574704938732	0.9557574704938733	3	1	62.503947	187.511841	0.7732971809418951	0.7732971809418953	3651.350033829367	32862.1503044643	0.7390845575970686	0.7390845575970688
FAPs	Efemp1	Egfr	Inflammatory-Mac	3	1	58.41791133333334	175.253734	0.9557574704938732	0.9557574704938733	3	1	0.4338690000000001	1.301607	0.005367815805265532	0.005367815805265533	25.345720772282	228.1114869505381	0.005130330056117617	0.005130330056117619
FAPs	Efemp1	Egfr	MuSCs	3	1	58.41791133333334	175.253734	0.9557574704938732	0.9557574704938733	3	1	15.972384	47.917152	0.1976099128607259	0.1976099128607259	933.0733122939521	8397.659810645569	0.1888671504602821	0.1888671504602821
FAPs	Efemp1	Egfr	Resolving-Mac	3	1	58.41791133333334	175.253734	0.9557574704938732	0.9557574704938733	3	1	0.2157183333333333	0.647155	0.002668861520763652	0.002668861520763652	12.60181446964111	113.41633022677	0.0025507843361835	0.0025507843361835
MuSCs	Efemp1	Egfr	ECs	3	1	0.2044496666666667	0.613349	0.003344938081433099	0.003344938081433099	3	1	1.701929666666667	5.105789	0.02105622887134972	0.02105622887134972	0.3479589530401111	3.131630577361	0.00007043178180314874	0.00007043178180314877
MuSCs	Efemp1	Egfr	FAPs	3	1	0.2044496666666667	0.613349	0.003344938081433099	0.003344938081433099	3	1	62.503947	187.511841	0.7732971809418951	0.7732971809418953	12.778911129501	115.010200165509	0.002586631188797406	0.002586631188797408
MuSCs	Efemp1	Egfr	Inflammatory-Mac	3	1	0.2044496666666667	0.613349	0.003344938081433099	0.003344938081433099	3	1	0.4338690000000001	1.301607	0.005367815805265532	0.005367815805265533	0.08870437242700001	0.7983393518430002	0.00001795501150115115	0.00001795501150115116
MuSCs	Efemp1	Egfr	MuSCs	3	1	0.2044496666666667	0.613349	0.003344938081433099	0.003344938081433099	3	1	15.972384	47.917152	0.1976099128607259	0.1976099128607259	3.265548584672	29.389937262048	0.0006609929227965183	0.0006609929227965186
MuSCs	Efemp1	Egfr	Resolving-Mac	3	1	0.2044496666666667	0.613349	0.003344938081433099	0.003344938081433099	3	1	0.2157183333333333	0.647155	0.002668861520763652	0.002668861520763652	0.04410354134388889	0.396931872095	0.000008927176534873794	0.000008927176534873795
Resolving-Mac	Efemp1	Egfr	ECs	2	0.6666666666666666	0.174281	0.5228429999999999	0.002851357809845171	0.002851357809845171	3	1	1.701929666666667	5.105789	0.02105622887134972	0.02105622887134972	0.2966140042363333	2.669526038127	0.00006003884263821038	0.00006003884263821038
Resolving-Mac	Efemp1	Egfr	FAPs	2	0.6666666666666666	0.174281	0.5228429999999999	0.002851357809845171	0.002851357809845171	3	1	62.503947	187.511841	0.7732971809418951	0.7732971809418953	10.893250387107	98.03925348396299	0.002204946956209927	0.002204946956209927
Resolving-Mac	Efemp1	Egfr	Inflammatory-Mac	2	0.6666666666666666	0.174281	0.5228429999999999	0.002851357809845171	0.002851357809845171	3	1	0.4338690000000001	1.301607	0.005367815805265532	0.005367815805265533	0.075615123189	0.680536108701	0.00001530556351815422	0.00001530556351815422
Resolving-Mac	Efemp1	Egfr	MuSCs	2	0.6666666666666666	0.174281	0.5228429999999999	0.002851357809845171	0.002851357809845171	3	1	15.972384	47.917152	0.1976099128607259	0.1976099128607259	2.783683055904	25.053147503136	0.0005634565683382544	0.0005634565683382545
Resolving-Mac	Efemp1	Egfr	Resolving-Mac	2	0.6666666666666666	0.174281	0.5228429999999999	0.002851357809845171	0.002851357809845171	3	1	0.2157183333333333	0.647155	0.002668861520763652	0.002668861520763652	0.03759560685166666	0.338360461665	0.000007609879140624699	0.000007609879140624699
"@

$lines = $data -split "`n"
$numRows = $lines.Length
$numCols = 20

$arr = New-Object 'object[,]' $numRows,$numCols

for ($r = 0; $r -lt $numRows; $r++) {
    $fields = $lines[$r].Trim() -split "`t"
    for ($c = 0; $c -lt $numCols; $c++) {
        $text = $fields[$c]
        if ($c -lt 4) {
            $arr[$r, $c] = $text
        } else {
            $arr[$r, $c] = [double]$text
        }
    }
}

$ws.Range("A2:T21").Value = $arr
